$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fill in the three empty answer cells in Yuki's standup row (table row 4).
# ---------------------------------------------------------------------------

# Cell (4,2) "What have you completed since last meeting?" — two runs, the
# first carrying an eastAsia font hint / zh-CN east-Asian language tag.
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 2)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Refined</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t xml:space="preserve"> eraser functionality so that instead of being a white pen, it is now functioning as an object eraser and erases every line caught in its path</w:t></w:r></w:p>'
$cell.Range.InsertXML($xml)

# Cell (4,3) "What are you going to complete today?"
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 3)
$cell.Range.Text = "Make eraser erase textbox and shapes"

# Cell (4,4) "What obstacles are in your way?"
$t = $d.Tables.Item(1)
$cell = $t.Cell(4, 4)
$cell.Range.Text = "The relation/translate position of the textbox and shapes need to be figured out clearly"

# ---------------------------------------------------------------------------
# 2) Simplify both "Edward need the JSONs..." paragraphs: drop the
#    gramStart/gramEnd proofErr markers around "need" and collapse the three
#    runs into a single clean run.
# ---------------------------------------------------------------------------

$targetText = "Edward need the JSONs and refactoring ready for him to start the object database."

$fixed = 0
while ($fixed -lt 2) {
    $found = $false
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($targetText + "`r")) {
            $rng = $p.Range
            $rng.End = $rng.End - 1
            # First write a distinct placeholder so the engine records a real
            # content change (identical-looking text is otherwise a no-op),
            # then write the final text so the paragraph collapses to one run.
            $rng.Text = "__TMP_PLACEHOLDER__"
            $rng2 = $p.Range
            $rng2.End = $rng2.End - 1
            $rng2.Text = $targetText
            $fixed = $fixed + 1
            $found = $true
            break
        }
    }
    if (-not $found) {
        break
    }
}

Write-Host "Fixed Edward paragraphs: " $fixed
